$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Senior Five")

# Row 3 - AMITO LUCY
$ws1.Range("D3").Value = 67.0
$ws1.Range("G3").Value = 44.0
$ws1.Range("G3").WrapText = $false
$ws1.Range("I3").Value = 71.0
$ws1.Range("I3").Font.ThemeColor = 1

# Row 7 - BONGOMIN RONNIE
$ws1.Range("D7").Value = 73.0
$ws1.Range("G7").Value = 52.0
$ws1.Range("G7").WrapText = $false
$ws1.Range("I7").Value = 60.0
$ws1.Range("I7").Font.ThemeColor = 1

# Row 9 - KIZITO STEPHEN
$ws1.Range("D9").Value = 85.0
$ws1.Range("G9").Value = 56.0
$ws1.Range("G9").WrapText = $false
$ws1.Range("I9").Value = 59.0
$ws1.Range("I9").Font.ThemeColor = 1

$ws2 = $wb.Worksheets.Item("Senior Six")

# Row 2 - ACII AGNESS
$ws2.Range("D2").Value = 42.0

# Row 7 - OCHEN ATIKU HUSSEIN
$ws2.Range("D7").Value = 33.0
